$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 previously contained the Norwegian/Danish note "Startet kl. 12:00"; it is
# corrected to a normal hours entry matching the rest of the column.
$ws.Range("B11").Value = "12 Hours"

# Row 12 was an empty placeholder row; it now becomes a new day entry: a date
# in column A (formatted like the other date cells) and a "Started at 00:00"
# note in column B (re-using the slot that used to hold "74 hours").
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A12").Value = (Get-Date -Year 2020 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B12").Value = "Started at 00:00"

# B13 used to hold "74 hours"; that text has moved, so the cell is cleared.
$ws.Range("B13").Value = $null

# Move the active selection from B14 to B13.
$ws.Range("B13").Select() | Out-Null
